$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Skills list ("Programming & Web: ..."): add Sass + Python, drop WebGL
#    and REST APIs, and drop the old trailing Sass entry.
#    Before: HTML5, CSS3, JavaScript, Node, React, Angular, Eleventy, WebGL,
#            C++, Java, Git, SQL, REST APIs, Sass, Jest, Chai
#    After:  HTML5, CSS3, Sass, JavaScript, Python, Node, React, Angular,
#            Eleventy, C++, Java, Git, SQL, Jest, Chai
# ---------------------------------------------------------------------------

# Insert "Sass, " right before "JavaScript" (after "CSS3, ")
$rng = $d.Content
$rng.Find.Execute("CSS3, ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Collapse(0)
$rng.InsertAfter("Sass, ")

# Insert "Python, " right before "Node," (after "JavaScript, ")
$rng = $d.Content
$rng.Find.Execute("JavaScript, ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Collapse(0)
$rng.InsertAfter("Python, ")

# Remove "WebGL, " (was right after "Eleventy, ")
$rng = $d.Content
$rng.Find.Execute("WebGL, ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Text = ""

# Remove "REST APIs, " and the old "Sass, " (were right after "SQL, ")
$rng = $d.Content
$rng.Find.Execute("REST APIs, Sass, ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Text = ""

# ---------------------------------------------------------------------------
# 2. Project blurb: "utilizing eleventy-plugin-sharp-respimg for" ->
#    "utilizing my eleventy-plugin-sharp-respimg plugin for"
# ---------------------------------------------------------------------------

$rng = $d.Content
$rng.Find.Execute("by utilizing ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Collapse(0)
$rng.InsertAfter("my ")

$rng = $d.Content
$rng.Find.Execute("respimg ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Collapse(0)
$rng.InsertAfter("plugin ")

# ---------------------------------------------------------------------------
# 3. Personal website blurb:
#    "using HTML, CSS, JavaScript" -> "using HTML, CSS/Sass, JavaScript"
#    "Used Lighthouse to verify" -> "Utilizing Lighthouse to verify"
# ---------------------------------------------------------------------------

$rng = $d.Content
$rng.Find.Execute("ground up using HTML, CSS", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Collapse(0)
$rng.InsertAfter("/Sass")

$rng = $d.Content
$rng.Find.Execute("Used Lighthouse", $true, $false, $false, $false, $false, $true, 1, $false, "Utilizing Lighthouse", 2) | Out-Null

# ---------------------------------------------------------------------------
# 4. Training blurb:
#    "training program to learn" -> "learning program studying"
#    "Built fullstack projects the MEAN stack" ->
#    "Built fullstack projects with the MEAN stack"
# ---------------------------------------------------------------------------

$rng = $d.Content
$rng.Find.Execute("training program to learn", $true, $false, $false, $false, $false, $true, 1, $false, "learning program studying", 2) | Out-Null

$rng = $d.Content
$rng.Find.Execute("projects ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Collapse(0)
$rng.InsertAfter("with ")
